$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Then_Goto" / "Else_Goto" columns (I/J) for the sex question (row 4)
# used to hold text labels ("pregnant" / "onset_date"); the template now
# points them at numeric row indexes of the rows to jump to.
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 5

# Leave the cursor on J6, matching the saved selection in the template.
$ws.Activate()
$ws.Range("J6").Select()
